$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D14's "Fecha Entrega" now uses the datetime number format (same style as
# the rest of the column) instead of the date-only format.
$ws.Cells.Item(14, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"

# New order rows (the search form now appends rows on every page, not just
# the last one) - rows 15 through 21 use the datetime format.
$rows = @(
    @(14, 1, 1, 45427, 0, 0, 0, 1),
    @(15, 1, 1, 45428, 0, 0, 0, 1),
    @(16, 5, 1, 45428, 0, 0, 0, 1),
    @(17, 1, 1, 45428, 0, 0, 0, 1),
    @(18, 3, 3, 45428, 0, 0, 0, 1),
    @(19, 6, 7, 45428, 0, 0, 0, 1),
    @(20, 6, 7, 45428, 0, 0, 0, 1),
    @(21, 6, 7, 45428, 0, 0, 0, 1)
)

$r = 15
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $ws.Cells.Item($r, 8).Value = $row[7]

    # "Fecha Recogida" (I) is left blank (empty text), matching the rest of
    # the sheet's un-collected orders.
    $ws.Cells.Item($r, 9).Value = "'"
    $ws.Cells.Item($r, 9).Style = "Normal"

    $r++
}

# Last new row (22): "Fecha Entrega" reverts to the date-only format, and
# "Pagado" is stored as a boolean TRUE instead of the numeric 1.
$ws.Cells.Item(22, 4).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(22, 8).Value = $true
$ws.Cells.Item(22, 9).Value = "'"
$ws.Cells.Item(22, 9).Style = "Normal"
